$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) body format used by columns B-E in data rows.
$plainStyle = $ws.Range("B2").Style

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = $plainStyle
}

# Row 2 - Bitcoin
Set-TextCell "D2" "58.701.34"
Set-TextCell "E2" "  +1.82%  "

# Row 3 - Ethereum
Set-TextCell "D3" "3.154.51"
Set-TextCell "E3" "  +1.91%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.01%  "

# Row 5 - BNB
Set-TextCell "D5" "534.14"
Set-TextCell "E5" "  +1.18%  "

# Row 6 - Solana
Set-TextCell "D6" "140.05"
Set-TextCell "E6" "  +1.80%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.02%  "

# Row 8 - XRP
Set-TextCell "D8" "0.517"
Set-TextCell "E8" "  +10.42%  "

# Row 9 - Toncoin
Set-TextCell "E9" "  +1.29%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +2.97%  "

# Row 11 - Cardano
Set-TextCell "D11" "0.422"
Set-TextCell "E11" "  +3.90%  "

# Row 12 - TRON
Set-TextCell "E12" "  +2.18%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "3.698.26"
Set-TextCell "E13" "  +2.32%  "

# Row 14 - Avalanche
Set-TextCell "D14" "25.84"
Set-TextCell "E14" "  +1.57%  "

# Row 15 - ShibaInu
Set-TextCell "E15" "  +5.54%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "58.747.06"
Set-TextCell "E16" "  +1.87%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "3.156.04"
Set-TextCell "E17" "  +2.43%  "

# Row 18 - Polkadot
Set-TextCell "E18" "  +4.71%  "

# Row 19 - Chainlink
Set-TextCell "D19" "13.01"
Set-TextCell "E19" "  +3.65%  "

# Row 20 - Uniswap
Set-TextCell "E20" "  +3.64%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "372.09"
Set-TextCell "E21" "  +6.11%  "

# Row 22 - LEO
Set-TextCell "D22" "5.81"
Set-TextCell "E22" "  +2.02%  "

# Row 23 - Dai
Set-TextCell "E23" "  +0.14%  "

# Row 24 - Litecoin
Set-TextCell "D24" "69.70"
Set-TextCell "E24" "  +1.79%  "

# Row 25 - Polygon
Set-TextCell "E25" "  +2.32%  "

# Row 26 - Kaspa
Set-TextCell "E26" "  +0.02%  "

# Row 27 - Binance-PegBSC-USD
Set-TextCell "E27" "  +0.05%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextCell "D28" "7.96"
Set-TextCell "E28" "  +11.30%  "

# Row 29 - PEPE
Set-TextCell "D29" "0.0₃0875"
Set-TextCell "E29" "  +1.02%  "

# Row 30 - was PancakeSwap, now RenderToken
Set-TextCell "B30" "RenderToken"
Set-TextCell "C30" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D30" "6.16"
Set-TextCell "E30" "  +2.62%  "

# Row 31 - was RenderToken, now PancakeSwap
Set-TextCell "B31" "PancakeSwap"
Set-TextCell "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D31" "1.89"
Set-TextCell "E31" "  +1.46%  "

# Row 32 - EthereumClassic
Set-TextCell "D32" "21.92"
Set-TextCell "E32" "  +3.51%  "

# Row 33 - NEARProtocol
Set-TextCell "D33" "5.20"
Set-TextCell "E33" "  +6.85%  "

# Row 34 - Fetch.AI
Set-TextCell "E34" "  +2.47%  "

# Row 35 - Monero
Set-TextCell "D35" "160.01"
Set-TextCell "E35" "  +0.42%  "

# Row 36 - Aptos
Set-TextCell "D36" "6.25"
Set-TextCell "E36" "  +3.59%  "

# Row 37 - ImmutableX
Set-TextCell "D37" "1.37"
Set-TextCell "E37" "  +10.01%  "

# Row 38 - EnergySwap
Set-TextCell "D38" "25.33"
Set-TextCell "E38" "  -1.69%  "

# Row 39 - Maker
Set-TextCell "D39" "2.653.42"
Set-TextCell "E39" "  +10.71%  "

# Row 40 - Stacks
Set-TextCell "D40" "1.67"
Set-TextCell "E40" "  +4.49%  "

# Row 41 - Hedera
Set-TextCell "E41" "  +2.96%  "

# Row 42 - Filecoin
Set-TextCell "E42" "  +3.70%  "

# Row 43 - was OKB, now Mantle
Set-TextCell "B43" "Mantle"
Set-TextCell "C43" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D43" "0.709"
Set-TextCell "E43" "  +2.34%  "

# Row 44 - was Mantle, now OKB
Set-TextCell "B44" "OKB"
Set-TextCell "C44" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D44" "38.56"
Set-TextCell "E44" "  +4.47%  "

# Row 45 - VeChain
Set-TextCell "D45" "0.0285"
Set-TextCell "E45" "  +7.42%  "

# Row 46 - FirstDigitalUSD
Set-TextCell "E46" "  -0.11%  "

# Row 47 - RenzoRestakedETH
Set-TextCell "D47" "3.196.45"
Set-TextCell "E47" "  +1.82%  "

# Row 48 - Stellar
Set-TextCell "D48" "0.104"
Set-TextCell "E48" "  +13.29%  "

# Row 49 - ONDO
Set-TextCell "D49" "0.985"
Set-TextCell "E49" "  +3.06%  "

# Row 50 - Cosmos
Set-TextCell "D50" "6.19"
Set-TextCell "E50" "  +2.76%  "

# Row 51 - InjectiveProtocol
Set-TextCell "D51" "20.26"
Set-TextCell "E51" "  +3.63%  "
